# Insert a new daily price record for Perejil (Vega Central Mapocho de Santiago)
# as row 186, pushing the existing rows 186:248 down to 187:249 and growing the
# used range from A1:R248 to A1:R249.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 186:248 down by one (EntireRow insert, like right-click > Insert in Excel).
$ws.Rows(186).Insert()

# Populate the newly-inserted row 186 with the new record's data.
$ws.Range("A186").Value = 9
$ws.Range("B186").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C186").Value = "Metropolitana"
$ws.Range("D186").Value = 44524
$ws.Range("E186").Value = 13
$ws.Range("F186").Value = 100112044
$ws.Range("G186").Value = "Perejil"
$ws.Range("H186").Value = "Sin especificar"
$ws.Range("I186").Value = "Primera"
$ws.Range("J186").Value = 97
$ws.Range("K186").Value = 11000
$ws.Range("L186").Value = 12000
$ws.Range("M186").Value = 11495
$ws.Range("N186").Value = "`$/docena de atados"
$ws.Range("O186").Value = "Región Metropolitana"
$ws.Range("P186").Value = 3832
$ws.Range("Q186").Value = 3
$ws.Range("R186").Value = "Hortaliza"
